# Cosmetic changes to the rubric
#
# The "Actual" score column on the Rubric sheet is cleared out (it was a
# filled-in grading pass with per-item scores + a couple of grader remarks in
# column D); the Grade sheet (which already holds the real/filled totals) is
# left alone content-wise. Only the active-sheet/selection bookkeeping and a
# little column-width tidy-up change on top of that.

$wb     = $excel.ActiveWorkbook
$rubric = $wb.Worksheets.Item("Rubric")
$grade  = $wb.Worksheets.Item("Grade")

# --- Rubric sheet: blank out the "Actual" column and the grader's remarks ---
$rubric.Range("C4").ClearContents()
$rubric.Range("C6:C15").ClearContents()
$rubric.Range("C17").ClearContents()
$rubric.Range("D6").ClearContents()
$rubric.Range("D14").ClearContents()
$rubric.Range("D15").ClearContents()

# --- cosmetic column-width tidy-up ---
$rubric.Columns.Item(1).ColumnWidth = 22.830729166666668
$rubric.Columns.Item(2).ColumnWidth = 8.498697916666666
$rubric.Columns.Item(3).ColumnWidth = 10.166666666666666
$rubric.Columns.Item(4).ColumnWidth = 9.498697916666666

$grade.Columns.Item(1).ColumnWidth = 22.330729166666668
$grade.Columns.Item(2).ColumnWidth = 7.998697916666667
$grade.Columns.Item(3).ColumnWidth = 5.998697916666667

# --- selection / active-sheet bookkeeping ---
[void]$rubric.Range("D20").Select()
[void]$grade.Activate()
[void]$grade.Range("E7").Select()
